$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "20-0=20"
$t.Cell(1,2).Range.Text = "17+26=43"
$t.Cell(1,3).Range.Text = "36+20=56"
$t.Cell(1,4).Range.Text = "53-16=37"
$t.Cell(1,5).Range.Text = "48+15=63"
$t.Cell(2,1).Range.Text = "36-2=34"
$t.Cell(2,2).Range.Text = "94-11=83"
$t.Cell(2,3).Range.Text = "96-42=54"
$t.Cell(2,4).Range.Text = "44+33=77"
$t.Cell(2,5).Range.Text = "47-16=31"
$t.Cell(3,1).Range.Text = "85-77=8"
$t.Cell(3,2).Range.Text = "1+79=80"
$t.Cell(3,3).Range.Text = "90-9=81"
$t.Cell(3,4).Range.Text = "7+20=27"
$t.Cell(3,5).Range.Text = "22+68=90"
$t.Cell(4,1).Range.Text = "57-30=27"
$t.Cell(4,2).Range.Text = "55+16=71"
$t.Cell(4,3).Range.Text = "29-28=1"
$t.Cell(4,4).Range.Text = "67+32=99"
$t.Cell(4,5).Range.Text = "26+55=81"
$t.Cell(5,1).Range.Text = "60-18=42"
$t.Cell(5,2).Range.Text = "7+74=81"
$t.Cell(5,3).Range.Text = "89-25=64"
$t.Cell(5,4).Range.Text = "10-8=2"
$t.Cell(5,5).Range.Text = "36+54=90"
$t.Cell(6,1).Range.Text = "5+17=22"
$t.Cell(6,2).Range.Text = "93+0=93"
$t.Cell(6,3).Range.Text = "99-24=75"
$t.Cell(6,4).Range.Text = "98-84=14"
$t.Cell(6,5).Range.Text = "60+31=91"
$t.Cell(7,1).Range.Text = "40-19=21"
$t.Cell(7,2).Range.Text = "31+63=94"
$t.Cell(7,3).Range.Text = "23+13=36"
$t.Cell(7,4).Range.Text = "14+70=84"
$t.Cell(7,5).Range.Text = "7+61=68"
$t.Cell(8,1).Range.Text = "88-76=12"
$t.Cell(8,2).Range.Text = "22+19=41"
$t.Cell(8,3).Range.Text = "59-13=46"
$t.Cell(8,4).Range.Text = "58+29=87"
$t.Cell(8,5).Range.Text = "69-32=37"
$t.Cell(9,1).Range.Text = "31+12=43"
$t.Cell(9,2).Range.Text = "13-12=1"
$t.Cell(9,3).Range.Text = "77-69=8"
$t.Cell(9,4).Range.Text = "2+89=91"
$t.Cell(9,5).Range.Text = "57+11=68"
$t.Cell(10,1).Range.Text = "50+32=82"
$t.Cell(10,2).Range.Text = "8+63=71"
$t.Cell(10,3).Range.Text = "55+39=94"
$t.Cell(10,4).Range.Text = "57-3=54"
$t.Cell(10,5).Range.Text = "12+30=42"
$t.Cell(11,1).Range.Text = "12-7=5"
$t.Cell(11,2).Range.Text = "89-32=57"
$t.Cell(11,3).Range.Text = "22+22=44"
$t.Cell(11,4).Range.Text = "55-30=25"
$t.Cell(11,5).Range.Text = "81-60=21"
$t.Cell(12,1).Range.Text = "88-20=68"
$t.Cell(12,2).Range.Text = "33+17=50"
$t.Cell(12,3).Range.Text = "93-54=39"
$t.Cell(12,4).Range.Text = "23+18=41"
$t.Cell(12,5).Range.Text = "60+34=94"
$t.Cell(13,1).Range.Text = "24+11=35"
$t.Cell(13,2).Range.Text = "96-32=64"
$t.Cell(13,3).Range.Text = "33+57=90"
$t.Cell(13,4).Range.Text = "95-38=57"
$t.Cell(13,5).Range.Text = "9+79=88"
$t.Cell(14,1).Range.Text = "86-70=16"
$t.Cell(14,2).Range.Text = "0+71=71"
$t.Cell(14,3).Range.Text = "24+3=27"
$t.Cell(14,4).Range.Text = "0+25=25"
$t.Cell(14,5).Range.Text = "27-18=9"
$t.Cell(15,1).Range.Text = "65-4=61"
$t.Cell(15,2).Range.Text = "99-54=45"
$t.Cell(15,3).Range.Text = "28+68=96"
$t.Cell(15,4).Range.Text = "33+32=65"
$t.Cell(15,5).Range.Text = "93-66=27"
$t.Cell(16,1).Range.Text = "14+12=26"
$t.Cell(16,2).Range.Text = "37-3=34"
$t.Cell(16,3).Range.Text = "66-19=47"
$t.Cell(16,4).Range.Text = "50-24=26"
$t.Cell(16,5).Range.Text = "13+72=85"
$t.Cell(17,1).Range.Text = "55+31=86"
$t.Cell(17,2).Range.Text = "70+25=95"
$t.Cell(17,3).Range.Text = "53+38=91"
$t.Cell(17,4).Range.Text = "78-27=51"
$t.Cell(17,5).Range.Text = "64-54=10"
$t.Cell(18,1).Range.Text = "9+69=78"
$t.Cell(18,2).Range.Text = "1+79=80"
$t.Cell(18,3).Range.Text = "14+17=31"
$t.Cell(18,4).Range.Text = "77-64=13"
$t.Cell(18,5).Range.Text = "45+40=85"
$t.Cell(19,1).Range.Text = "93-79=14"
$t.Cell(19,2).Range.Text = "6+41=47"
$t.Cell(19,3).Range.Text = "63-34=29"
$t.Cell(19,4).Range.Text = "8+21=29"
$t.Cell(19,5).Range.Text = "77-76=1"
$t.Cell(20,1).Range.Text = "54+18=72"
$t.Cell(20,2).Range.Text = "98-34=64"
$t.Cell(20,3).Range.Text = "15+29=44"
$t.Cell(20,4).Range.Text = "77-33=44"
$t.Cell(20,5).Range.Text = "64+32=96"
